$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 5.978421000000001
$ws.Range("N2").Value = 17.935263
$ws.Range("O2").Value = 0.05704457007880161
$ws.Range("P2").Value = 0.06242884486533885
$ws.Range("Q2").Value = 2.826226786698001
$ws.Range("R2").Value = 25.43604108028201
$ws.Range("S2").Value = 0.001868579866130652
$ws.Range("T2").Value = 0.002120794514267882
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.6646576013185088
$ws.Range("P3").Value = 0.7273927426214574
$ws.Range("Q3").Value = 32.92991978437001
$ws.Range("R3").Value = 296.36927805933
$ws.Range("S3").Value = 0.02177184980058229
$ws.Range("T3").Value = 0.02471054112241555
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("M4").Value = 1.290243
$ws.Range("N4").Value = 3.870729
$ws.Range("O4").Value = 0.01231116999491725
$ws.Range("P4").Value = 0.01347318632889677
$ws.Range("Q4").Value = 0.609946895334
$ws.Range("R4").Value = 5.489522058006
$ws.Range("S4").Value = 0.0004032707118177208
$ws.Range("T4").Value = 0.0004577028410131262
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 27.1166075
$ws.Range("N5").Value = 54.233215
$ws.Range("O5").Value = 0.2587397603536297
$ws.Range("P5").Value = 0.1887743138075849
$ws.Range("Q5").Value = 12.819050796335
$ws.Range("R5").Value = 76.91430477801001
$ws.Range("S5").Value = 0.008475406267351767
$ws.Range("T5").Value = 0.006412925467728608
$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 0.759494
$ws.Range("N6").Value = 2.278482
$ws.Range("O6").Value = 0.00724689825414258
$ws.Range("P6").Value = 0.007930912376722157
$ws.Range("Q6").Value = 0.359041674572
$ws.Range("R6").Value = 3.231375071148
$ws.Range("S6").Value = 0.0002373829472442695
$ws.Range("T6").Value = 0.0002694241019191138
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 5.978421000000001
$ws.Range("N7").Value = 17.935263
$ws.Range("O7").Value = 0.05704457007880161
$ws.Range("P7").Value = 0.06242884486533885
$ws.Range("Q7").Value = 70.93763192988001
$ws.Range("R7").Value = 638.4386873689201
$ws.Range("S7").Value = 0.04690091800100283
$ws.Range("T7").Value = 0.05323144673319507
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.6646576013185088
$ws.Range("P8").Value = 0.7273927426214574
$ws.Range("S8").Value = 0.5464683424753667
$ws.Range("T8").Value = 0.620228807957083
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("M9").Value = 1.290243
$ws.Range("N9").Value = 3.870729
$ws.Range("O9").Value = 0.01231116999491725
$ws.Range("P9").Value = 0.01347318632889677
$ws.Range("Q9").Value = 15.30952454404
$ws.Range("R9").Value = 137.78572089636
$ws.Range("S9").Value = 0.01012200063267004
$ws.Range("T9").Value = 0.01148823435609132
$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 27.1166075
$ws.Range("N10").Value = 54.233215
$ws.Range("O10").Value = 0.2587397603536297
$ws.Range("P10").Value = 0.1887743138075849
$ws.Range("Q10").Value = 321.7551795067667
$ws.Range("R10").Value = 1930.5310770406
$ws.Range("S10").Value = 0.2127307168268809
$ws.Range("T10").Value = 0.1609629307048587
$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 0.759494
$ws.Range("N11").Value = 2.278482
$ws.Range("O11").Value = 0.00724689825414258
$ws.Range("P11").Value = 0.007930912376722157
$ws.Range("Q11").Value = 9.011862132986668
$ws.Range("R11").Value = 81.10675919687999
$ws.Range("S11").Value = 0.005958256505564535
$ws.Range("T11").Value = 0.006762482000712437
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 5.978421000000001
$ws.Range("N12").Value = 17.935263
$ws.Range("O12").Value = 0.05704457007880161
$ws.Range("P12").Value = 0.06242884486533885
$ws.Range("Q12").Value = 2.26534328532
$ws.Range("R12").Value = 20.38808956788
$ws.Range("S12").Value = 0.001497747764880814
$ws.Range("T12").Value = 0.001699908738765206
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.6646576013185088
$ws.Range("P13").Value = 0.7273927426214574
$ws.Range("Q13").Value = 26.3947582058
$ws.Range("R13").Value = 237.5528238522
$ws.Range("S13").Value = 0.01745108141599922
$ws.Range("T13").Value = 0.01980656990152198
$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("M14").Value = 1.290243
$ws.Range("N14").Value = 3.870729
$ws.Range("O14").Value = 0.01231116999491725
$ws.Range("P14").Value = 0.01347318632889677
$ws.Range("Q14").Value = 0.48889887756
$ws.Range("R14").Value = 4.40008989804
$ws.Range("S14").Value = 0.0003232389571432183
$ws.Range("T14").Value = 0.0003668686683039945
$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 27.1166075
$ws.Range("N15").Value = 54.233215
$ws.Range("O15").Value = 0.2587397603536297
$ws.Range("P15").Value = 0.1887743138075849
$ws.Range("Q15").Value = 10.2750249139
$ws.Range("R15").Value = 61.6501494834
$ws.Range("S15").Value = 0.006793405528696511
$ws.Range("T15").Value = 0.00514023776009486
$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 0.759494
$ws.Range("N16").Value = 2.278482
$ws.Range("O16").Value = 0.00724689825414258
$ws.Range("P16").Value = 0.007930912376722157
$ws.Range("Q16").Value = 0.28778746648
$ws.Range("R16").Value = 2.59008719832
$ws.Range("S16").Value = 0.0001902727226704825
$ws.Range("T16").Value = 0.0002159550971133918
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 5.978421000000001
$ws.Range("N17").Value = 17.935263
$ws.Range("O17").Value = 0.05704457007880161
$ws.Range("P17").Value = 0.06242884486533885
$ws.Range("Q17").Value = 9.256730004297003
$ws.Range("R17").Value = 55.54038002578201
$ws.Range("S17").Value = 0.006120152633768509
$ws.Range("T17").Value = 0.004630820217158025
$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.6646576013185088
$ws.Range("P18").Value = 0.7273927426214574
$ws.Range("Q18").Value = 107.855242877805
$ws.Range("R18").Value = 647.13145726683
$ws.Range("S18").Value = 0.07130925807039035
$ws.Range("T18").Value = 0.05395622849679946
$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("M19").Value = 1.290243
$ws.Range("N19").Value = 3.870729
$ws.Range("O19").Value = 0.01231116999491725
$ws.Range("P19").Value = 0.01347318632889677
$ws.Range("Q19").Value = 1.997756780751
$ws.Range("R19").Value = 11.986540684506
$ws.Range("S19").Value = 0.001320831051317962
$ws.Range("T19").Value = 0.0009994082667390972
$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 27.1166075
$ws.Range("N20").Value = 54.233215
$ws.Range("O20").Value = 0.2587397603536297
$ws.Range("P20").Value = 0.1887743138075849
$ws.Range("Q20").Value = 41.98618903887751
$ws.Range("R20").Value = 167.94475615551
$ws.Range("S20").Value = 0.02775946638920074
$ws.Range("T20").Value = 0.01400282050302122
$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 0.759494
$ws.Range("N21").Value = 2.278482
$ws.Range("O21").Value = 0.00724689825414258
$ws.Range("P21").Value = 0.007930912376722157
$ws.Range("Q21").Value = 1.175967851358
$ws.Range("R21").Value = 7.055807108148
$ws.Range("S21").Value = 0.0007774994776097868
$ws.Range("T21").Value = 0.0005882958343031071
$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 5.978421000000001
$ws.Range("N22").Value = 17.935263
$ws.Range("O22").Value = 0.05704457007880161
$ws.Range("P22").Value = 0.06242884486533885
$ws.Range("Q22").Value = 0.9939722754600001
$ws.Range("R22").Value = 8.945750479140001
$ws.Range("S22").Value = 0.0006571718130188014
$ws.Range("T22").Value = 0.000745874661952663
$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.6646576013185088
$ws.Range("P23").Value = 0.7273927426214574
$ws.Range("Q23").Value = 11.5813166349
$ws.Range("R23").Value = 104.2318497141
$ws.Range("S23").Value = 0.007657069556170249
$ws.Range("T23").Value = 0.008690595143637297
$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("M24").Value = 1.290243
$ws.Range("N24").Value = 3.870729
$ws.Range("O24").Value = 0.01231116999491725
$ws.Range("P24").Value = 0.01347318632889677
$ws.Range("Q24").Value = 0.21451580118
$ws.Range("R24").Value = 1.93064221062
$ws.Range("S24").Value = 0.0001418286419683085
$ws.Range("T24").Value = 0.0001609721967492403
$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 27.1166075
$ws.Range("N25").Value = 54.233215
$ws.Range("O25").Value = 0.2587397603536297
$ws.Range("P25").Value = 0.1887743138075849
$ws.Range("Q25").Value = 4.50840716295
$ws.Range("R25").Value = 27.0504429777
$ws.Range("S25").Value = 0.002980765341499741
$ws.Range("T25").Value = 0.002255399371881588
$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 0.759494
$ws.Range("N26").Value = 2.278482
$ws.Range("O26").Value = 0.00724689825414258
$ws.Range("P26").Value = 0.007930912376722157
$ws.Range("Q26").Value = 0.12627347244
$ws.Range("R26").Value = 1.13646125196
$ws.Range("S26").Value = 0.00008348660105350584
$ws.Range("T26").Value = 0.00009475534267410672
